$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "NA" for the duplicate_image_filename column (E) on rows 2-21,
# matching the rest of the table which already has data in that range.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
